$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted as row 11 ("Ramas de apio",
# week of 2023-07-21), pushing the existing rows 11-20 down to rows 12-21.
$ws.Rows("11:11").Insert()

$ws.Cells.Item(11, 1).Value  = 1
$ws.Cells.Item(11, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(11, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(11, 4).Value  = 45128
$ws.Cells.Item(11, 5).Value  = 15
$ws.Cells.Item(11, 6).Value  = 100112017
$ws.Cells.Item(11, 7).Value  = "Ramas de apio"
$ws.Cells.Item(11, 8).Value  = "Sin especificar"
$ws.Cells.Item(11, 9).Value  = "Primera"
$ws.Cells.Item(11, 10).Value = 200
$ws.Cells.Item(11, 11).Value = 3500
$ws.Cells.Item(11, 12).Value = 4000
$ws.Cells.Item(11, 13).Value = 3750
$ws.Cells.Item(11, 14).Value = "$/atado 7 kilos"
$ws.Cells.Item(11, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(11, 16).Value = 3750
$ws.Cells.Item(11, 17).Value = 1
$ws.Cells.Item(11, 18).Value = "Hortaliza"
